$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.919.75'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '3.523.39'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.07'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.31'
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('D7').Value = '3.521.25'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.123'
$ws.Range('E10').Value = '  +1.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.14'
$ws.Range('E11').Value = '  +3.60%  '
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').Value = '4.122.78'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.40'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000182'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '3.519.80'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').Value = '64.946.96'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.01'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '390.43'
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.575'
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').Value = '3.667.58'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.15'
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.62'
$ws.Range('E28').Value = '  +20.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.73'
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('E32').Value = '  +2.59%  '
$ws.Range('D33').Value = '3.529.31'
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.14'
$ws.Range('E34').Value = '  +2.19%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.22'
$ws.Range('E37').Value = '  +6.18%  '
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '168.54'
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.84'
$ws.Range('E40').Value = '  +0.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0819'
$ws.Range('E41').Value = '  +3.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.821'
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.67'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('E44').Value = '  +4.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.33'
$ws.Range('E46').Value = '  -4.18%  '
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('D50').Value = '2.404.12'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.897'
$ws.Range('E51').Value = '  +6.29%  '
